# Update two odds values in row 5 (S5, T5)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("S5").Value = 1.44
$ws.Range("T5").Value = 2.63

# Row 9 (PORTUGAL - LIGA PORTUGAL / Sporting CP vs Estrela) was removed.
# Deleting it shifts the old row 10 (Uruguay) up to row 9 and the old
# row 11 (USA - MLS) up to row 10, matching the target workbook exactly.
$ws.Rows.Item(9).Delete()
